# --- Reproduce the PanelApp "Hereditary ataxia - adult onset" commit ---
# 1) Add a new "metadata" worksheet after "data"
# 2) Refresh the F-column "time_taken" timestamps on the "data" sheet
# 3) Populate "metadata" with the panel/query metadata row

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Step 1: updated query timestamps for every gene row in "data" ---
$newTimestamps = @(
  "2021-10-05 14:20:42.226223",
  "2021-10-05 14:20:42.226231",
  "2021-10-05 14:20:42.226234",
  "2021-10-05 14:20:42.226237",
  "2021-10-05 14:20:42.226240",
  "2021-10-05 14:20:42.226242",
  "2021-10-05 14:20:42.226245",
  "2021-10-05 14:20:42.226247",
  "2021-10-05 14:20:42.226250",
  "2021-10-05 14:20:42.226253",
  "2021-10-05 14:20:42.226255",
  "2021-10-05 14:20:42.226258",
  "2021-10-05 14:20:42.226260",
  "2021-10-05 14:20:42.226262",
  "2021-10-05 14:20:42.226265",
  "2021-10-05 14:20:42.226267",
  "2021-10-05 14:20:42.226270",
  "2021-10-05 14:20:42.226273",
  "2021-10-05 14:20:42.226276",
  "2021-10-05 14:20:42.226278",
  "2021-10-05 14:20:42.226280",
  "2021-10-05 14:20:42.226283",
  "2021-10-05 14:20:42.226285",
  "2021-10-05 14:20:42.226288",
  "2021-10-05 14:20:42.226290",
  "2021-10-05 14:20:42.226293",
  "2021-10-05 14:20:42.226296",
  "2021-10-05 14:20:42.226298",
  "2021-10-05 14:20:42.226300",
  "2021-10-05 14:20:42.226303",
  "2021-10-05 14:20:42.226305",
  "2021-10-05 14:20:42.226308",
  "2021-10-05 14:20:42.226310",
  "2021-10-05 14:20:42.226313",
  "2021-10-05 14:20:42.226315",
  "2021-10-05 14:20:42.226318",
  "2021-10-05 14:20:42.226320",
  "2021-10-05 14:20:42.226323",
  "2021-10-05 14:20:42.226325",
  "2021-10-05 14:20:42.226327",
  "2021-10-05 14:20:42.226330",
  "2021-10-05 14:20:42.226333",
  "2021-10-05 14:20:42.226335",
  "2021-10-05 14:20:42.226338",
  "2021-10-05 14:20:42.226340",
  "2021-10-05 14:20:42.226343",
  "2021-10-05 14:20:42.226345",
  "2021-10-05 14:20:42.226348",
  "2021-10-05 14:20:42.226350",
  "2021-10-05 14:20:42.226352",
  "2021-10-05 14:20:42.226355",
  "2021-10-05 14:20:42.226357",
  "2021-10-05 14:20:42.226360",
  "2021-10-05 14:20:42.226363",
  "2021-10-05 14:20:42.226365",
  "2021-10-05 14:20:42.226367",
  "2021-10-05 14:20:42.226370",
  "2021-10-05 14:20:42.226372",
  "2021-10-05 14:20:42.226375",
  "2021-10-05 14:20:42.226377",
  "2021-10-05 14:20:42.226380",
  "2021-10-05 14:20:42.226382",
  "2021-10-05 14:20:42.226384",
  "2021-10-05 14:20:42.226387",
  "2021-10-05 14:20:42.226391",
  "2021-10-05 14:20:42.226393",
  "2021-10-05 14:20:42.226396",
  "2021-10-05 14:20:42.226398",
  "2021-10-05 14:20:42.226401",
  "2021-10-05 14:20:42.226403",
  "2021-10-05 14:20:42.226405",
  "2021-10-05 14:20:42.226408",
  "2021-10-05 14:20:42.226410",
  "2021-10-05 14:20:42.226413",
  "2021-10-05 14:20:42.226415",
  "2021-10-05 14:20:42.226418",
  "2021-10-05 14:20:42.226422",
  "2021-10-05 14:20:42.226425",
  "2021-10-05 14:20:42.226428",
  "2021-10-05 14:20:42.226430",
  "2021-10-05 14:20:42.226433",
  "2021-10-05 14:20:42.226435",
  "2021-10-05 14:20:42.226437",
  "2021-10-05 14:20:42.226440",
  "2021-10-05 14:20:42.226442",
  "2021-10-05 14:20:42.226445",
  "2021-10-05 14:20:42.226447",
  "2021-10-05 14:20:42.226450",
  "2021-10-05 14:20:42.226452",
  "2021-10-05 14:20:42.226455",
  "2021-10-05 14:20:42.226457",
  "2021-10-05 14:20:42.226460",
  "2021-10-05 14:20:42.226463",
  "2021-10-05 14:20:42.226466",
  "2021-10-05 14:20:42.226468",
  "2021-10-05 14:20:42.226471",
  "2021-10-05 14:20:42.226473",
  "2021-10-05 14:20:42.226476",
  "2021-10-05 14:20:42.226478",
  "2021-10-05 14:20:42.226481",
  "2021-10-05 14:20:42.226483",
  "2021-10-05 14:20:42.226486",
  "2021-10-05 14:20:42.226488",
  "2021-10-05 14:20:42.226491",
  "2021-10-05 14:20:42.226493",
  "2021-10-05 14:20:42.226495",
  "2021-10-05 14:20:42.226498",
  "2021-10-05 14:20:42.226500",
  "2021-10-05 14:20:42.226504",
  "2021-10-05 14:20:42.226507",
  "2021-10-05 14:20:42.226510",
  "2021-10-05 14:20:42.226512",
  "2021-10-05 14:20:42.226515",
  "2021-10-05 14:20:42.226517",
  "2021-10-05 14:20:42.226520",
  "2021-10-05 14:20:42.226522",
  "2021-10-05 14:20:42.226524",
  "2021-10-05 14:20:42.226527",
  "2021-10-05 14:20:42.226529",
  "2021-10-05 14:20:42.226532",
  "2021-10-05 14:20:42.226534",
  "2021-10-05 14:20:42.226537",
  "2021-10-05 14:20:42.226539",
  "2021-10-05 14:20:42.226542",
  "2021-10-05 14:20:42.226544",
  "2021-10-05 14:20:42.226546",
  "2021-10-05 14:20:42.226549",
  "2021-10-05 14:20:42.226551",
  "2021-10-05 14:20:42.226556",
  "2021-10-05 14:20:42.226559",
  "2021-10-05 14:20:42.226561",
  "2021-10-05 14:20:42.226564",
  "2021-10-05 14:20:42.226566",
  "2021-10-05 14:20:42.226569",
  "2021-10-05 14:20:42.226571",
  "2021-10-05 14:20:42.226573",
  "2021-10-05 14:20:42.226576",
  "2021-10-05 14:20:42.226578",
  "2021-10-05 14:20:42.226581",
  "2021-10-05 14:20:42.226583",
  "2021-10-05 14:20:42.226586",
  "2021-10-05 14:20:42.226588",
  "2021-10-05 14:20:42.226590",
  "2021-10-05 14:20:42.226593",
  "2021-10-05 14:20:42.226595",
  "2021-10-05 14:20:42.226598",
  "2021-10-05 14:20:42.226600",
  "2021-10-05 14:20:42.226603",
  "2021-10-05 14:20:42.226605",
  "2021-10-05 14:20:42.226608",
  "2021-10-05 14:20:42.226611",
  "2021-10-05 14:20:42.226613",
  "2021-10-05 14:20:42.226615",
  "2021-10-05 14:20:42.226618",
  "2021-10-05 14:20:42.226620",
  "2021-10-05 14:20:42.226623",
  "2021-10-05 14:20:42.226625",
  "2021-10-05 14:20:42.226628",
  "2021-10-05 14:20:42.226630",
  "2021-10-05 14:20:42.226633",
  "2021-10-05 14:20:42.226635",
  "2021-10-05 14:20:42.226638",
  "2021-10-05 14:20:42.226640",
  "2021-10-05 14:20:42.226642",
  "2021-10-05 14:20:42.226645",
  "2021-10-05 14:20:42.226647",
  "2021-10-05 14:20:42.226650",
  "2021-10-05 14:20:42.226652",
  "2021-10-05 14:20:42.226655",
  "2021-10-05 14:20:42.226657",
  "2021-10-05 14:20:42.226659",
  "2021-10-05 14:20:42.226662",
  "2021-10-05 14:20:42.226666",
  "2021-10-05 14:20:42.226669",
  "2021-10-05 14:20:42.226671",
  "2021-10-05 14:20:42.226674",
  "2021-10-05 14:20:42.226676",
  "2021-10-05 14:20:42.226678",
  "2021-10-05 14:20:42.226681",
  "2021-10-05 14:20:42.226683",
  "2021-10-05 14:20:42.226686",
  "2021-10-05 14:20:42.226688",
  "2021-10-05 14:20:42.226691",
  "2021-10-05 14:20:42.226693",
  "2021-10-05 14:20:42.226695",
  "2021-10-05 14:20:42.226698",
  "2021-10-05 14:20:42.226700",
  "2021-10-05 14:20:42.226703",
  "2021-10-05 14:20:42.226705",
  "2021-10-05 14:20:42.226708",
  "2021-10-05 14:20:42.226710",
  "2021-10-05 14:20:42.226713",
  "2021-10-05 14:20:42.226715",
  "2021-10-05 14:20:42.226718",
  "2021-10-05 14:20:42.226720",
  "2021-10-05 14:20:42.226723",
  "2021-10-05 14:20:42.226725",
  "2021-10-05 14:20:42.226727",
  "2021-10-05 14:20:42.226730",
  "2021-10-05 14:20:42.226732",
  "2021-10-05 14:20:42.226735",
  "2021-10-05 14:20:42.226738",
  "2021-10-05 14:20:42.226740",
  "2021-10-05 14:20:42.226743",
  "2021-10-05 14:20:42.226745",
  "2021-10-05 14:20:42.226747",
  "2021-10-05 14:20:42.226750",
  "2021-10-05 14:20:42.226752",
  "2021-10-05 14:20:42.226755",
  "2021-10-05 14:20:42.226757",
  "2021-10-05 14:20:42.226760",
  "2021-10-05 14:20:42.226762",
  "2021-10-05 14:20:42.226765",
  "2021-10-05 14:20:42.226767",
  "2021-10-05 14:20:42.226770",
  "2021-10-05 14:20:42.226772",
  "2021-10-05 14:20:42.226775",
  "2021-10-05 14:20:42.226777",
  "2021-10-05 14:20:42.226780",
  "2021-10-05 14:20:42.226782",
  "2021-10-05 14:20:42.226784",
  "2021-10-05 14:20:42.226787",
  "2021-10-05 14:20:42.226789",
  "2021-10-05 14:20:42.226792",
  "2021-10-05 14:20:42.226794",
  "2021-10-05 14:20:42.226797",
  "2021-10-05 14:20:42.226799",
  "2021-10-05 14:20:42.226802",
  "2021-10-05 14:20:42.226804",
  "2021-10-05 14:20:42.226806",
  "2021-10-05 14:20:42.226809",
  "2021-10-05 14:20:42.226811",
  "2021-10-05 14:20:42.226815",
  "2021-10-05 14:20:42.226818",
  "2021-10-05 14:20:42.226820",
  "2021-10-05 14:20:42.226823",
  "2021-10-05 14:20:42.226825",
  "2021-10-05 14:20:42.226828",
  "2021-10-05 14:20:42.226830",
  "2021-10-05 14:20:42.226833",
  "2021-10-05 14:20:42.226835",
  "2021-10-05 14:20:42.226838",
  "2021-10-05 14:20:42.226840"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
  $row = $i + 2
  $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- Step 2: add the "metadata" sheet right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 1).Font.Bold = $true
$metaSheet.Cells.Item(2, 1).Borders.LineStyle = 1
$metaSheet.Cells.Item(2, 1).HorizontalAlignment = -4108
$metaSheet.Cells.Item(2, 1).VerticalAlignment = -4160

$metaSheet.Cells.Item(2, 2).Value = "Hereditary ataxia - adult onset"
$metaSheet.Cells.Item(2, 3).Value = 466
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "2.85"
$metaSheet.Cells.Item(2, 5).Value = "2021-08-26T09:09:07.266393Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:20:42.222948"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/466/?format=json"

$metaSheet.Range("A1").Select()
